# Add a new derived-variable row ("Ca21" / cancer_active / Active cancer...)
# to the table on Sheet1, directly above the existing "D01" (age) row 51,
# shifting every following row down by one — matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Physically insert a new row at row 51 (pushes rows 51..183 down to 52..184).
$ws.Rows("51:51").Insert()

# Populate the new row with the new variable's data.
$ws.Range("A51").Value2 = "Ca21"
$ws.Range("B51").Value2 = "cancer_active"
$ws.Range("C51").Value2 = "Cancer"
$ws.Range("D51").Value2 = "Active cancer/treatment vs inactive"

# Grow the Table1 structured range (and its AutoFilter) to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E184"))

# Keep the existing selection state (D51 stays selected, same as before/after
# in the source edit).
$ws.Range("D51").Select()
